# Update the division-problem worksheet numbers for the regenerated table.
# Each entry is addressed by its (row, column) position inside the single
# table on the page so that duplicate cell text (e.g. "83÷6=" occurring
# twice) is updated independently without cross-matching.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; Old = "87÷8="; New = "65÷8=" },
    @{ Row = 1;  Col = 2; Old = "85÷7="; New = "81÷2=" },
    @{ Row = 1;  Col = 3; Old = "91÷6="; New = "78÷3=" },
    @{ Row = 1;  Col = 4; Old = "51÷6="; New = "52÷5=" },
    @{ Row = 1;  Col = 5; Old = "11÷7="; New = "54÷6=" },

    @{ Row = 5;  Col = 1; Old = "23÷6="; New = "81÷3=" },
    @{ Row = 5;  Col = 2; Old = "77÷5="; New = "59÷7=" },
    @{ Row = 5;  Col = 3; Old = "80÷7="; New = "81÷2=" },
    @{ Row = 5;  Col = 4; Old = "16÷2="; New = "87÷5=" },
    @{ Row = 5;  Col = 5; Old = "88÷6="; New = "17÷3=" },

    @{ Row = 9;  Col = 1; Old = "49÷2="; New = "91÷7=" },
    @{ Row = 9;  Col = 2; Old = "94÷4="; New = "96÷4=" },
    @{ Row = 9;  Col = 3; Old = "17÷7="; New = "20÷2=" },
    @{ Row = 9;  Col = 4; Old = "65÷3="; New = "75÷4=" },
    @{ Row = 9;  Col = 5; Old = "15÷7="; New = "59÷9=" },

    @{ Row = 13; Col = 1; Old = "83÷6="; New = "49÷4=" },
    @{ Row = 13; Col = 2; Old = "43÷3="; New = "89÷4=" },
    @{ Row = 13; Col = 3; Old = "15÷9="; New = "12÷8=" },
    @{ Row = 13; Col = 4; Old = "87÷7="; New = "70÷2=" },
    @{ Row = 13; Col = 5; Old = "48÷9="; New = "47÷9=" },

    @{ Row = 17; Col = 1; Old = "36÷5="; New = "80÷9=" },
    @{ Row = 17; Col = 2; Old = "54÷4="; New = "17÷3=" },
    @{ Row = 17; Col = 3; Old = "88÷3="; New = "34÷3=" },
    @{ Row = 17; Col = 4; Old = "84÷8="; New = "48÷3=" },
    @{ Row = 17; Col = 5; Old = "83÷6="; New = "22÷4=" }
)

foreach ($edit in $edits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $rng = $cell.Range

    # Defensive check: the cell's Range.Text includes the trailing
    # paragraph mark (chr 13) and end-of-cell mark (chr 7); strip them
    # before comparing against the expected original value.
    $current = $rng.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $edit.Old) {
        throw ("Unexpected cell text at row " + $edit.Row + ", col " + $edit.Col + `
               ": expected '" + $edit.Old + "' but found '" + $current + "'")
    }

    $rng.Text = $edit.New
}
